$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of data to append (Date serial, DeathCovid, DeathWithCovid, Total)
$newRows = @(
    @(44341, 12312, 2410, 14722),
    @(44342, 12320, 2413, 14733),
    @(44343, 12333, 2415, 14748),
    @(44344, 12335, 2416, 14751),
    @(44345, 12339, 2416, 14755),
    @(44346, 12343, 2417, 14760),
    @(44347, 12353, 2417, 14770),
    @(44348, 12366, 2417, 14783),
    @(44349, 12375, 2418, 14793),
    @(44350, 12382, 2420, 14802),
    @(44351, 12395, 2423, 14818),
    @(44352, 12404, 2423, 14827)
)

$startRow = 222
$lastExistingRow = $startRow - 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Copy the date-column formatting from the last existing data row so the
    # new cells reuse the same style instead of creating a brand new one.
    $ws.Cells.Item($lastExistingRow, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
$excel.CutCopyMode = $false

# Update the view: scroll so row 187 is at the top and select the last new cell
$excel.ActiveWindow.ScrollRow = 187
$ws.Range("A233").Select()
